# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# The underlying "Periodo Mora" database feeding this account-statement
# (Estado de Cuenta) sheet was refreshed: the six period codes listed in
# column E (rows 16-21) are re-pulled from the updated source in the
# opposite order, so the values effectively reverse top-to-bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "1807"
$ws.Range("E17").Value = "1806"
$ws.Range("E18").Value = "1805"
$ws.Range("E19").Value = "1804"
$ws.Range("E20").Value = "1803"
$ws.Range("E21").Value = "1802"
